$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows (and their styles)
# shift down by one automatically.
$ws.Rows.Item(1).Insert()

# Populate the new header row with plain (unstyled) string values.
$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "treated"
$ws.Range("C1").Value = "control"
$ws.Range("D1").Value = "pvalue"
$ws.Range("E1").Value = "treated_pre"
